# Finished Work break down structure (for real)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7: reassign task from Damien to Jake, adjust planned time
$ws.Range("D7").Value = "Jake"
$ws.Range("F7").Value = 1.5

# Rows 9-15: fill in planned time (F) for tasks that already had an owner
$ws.Range("F9").Value = 5
$ws.Range("F10").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 0.5
$ws.Range("F13").Value = 3
$ws.Range("F14").Value = 0.5
$ws.Range("F15").Value = 0.5

# Rows 22-27: assign owners and planned time to the remaining tasks
$ws.Range("D22").Value = "Kevin"
$ws.Range("F22").Value = 2

$ws.Range("D23").Value = "Kevin"
$ws.Range("F23").Value = 0.5

$ws.Range("D24").Value = "Damien"
$ws.Range("F24").Value = 0.5

$ws.Range("D25").Value = "Damien"
$ws.Range("F25").Value = 0.5

$ws.Range("D26").Value = "All"
$ws.Range("F26").Value = 3

$ws.Range("D27").Value = "Jake"
$ws.Range("F27").Value = 0.5

# Extend the planned-time subtotal to cover the newly filled rows
$ws.Range("F28").Formula = "=SUBTOTAL(9,F4:F27)"

# Recalculate so cached formula results stay correct
$excel.Calculate()

# Update the saved selection to reflect where the user ended up
$ws.Range("G28").Select()
